# Apply commit "Add descriptions titles" to the DMI Reference Organisation Interne workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Title (row with "Title" in column A) gets a value in column B
$meta.Range("B5").Value = "DMI Reference Organisation Interne"

# Date value is refreshed to the new generation timestamp
$meta.Range("B8").Value = "2026-02-25T08:15:31+00:00"

# Description (row with "Description" in column A) gets a value in column B
$meta.Range("B12").Value = "Extension créée dans ce volet pour référencer l'organisation fabricant du DMI."

# --- Sheet "Elements" ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 corresponds to the root "Extension" element definition.
$elements.Range("L2").Value = "DMI Reference Organisation Interne"
$elements.Range("M2").Value = "Extension créée dans ce volet pour référencer l'organisation fabricant du DMI."
$elements.Range("AK2").Value = ""
